# Weekly update for "Hortaliza, Terminal Hortofrutícola Agro Chillán - Apio".
#
# A new weekly record (fecha serial 44614, i.e. 2022-02-22) is inserted as
# row 59, pushing the existing rows 59..177 down to 60..178 (dimension grows
# from A1:R177 to A1:R178). Every column on the existing rows keeps its
# original value - only the position moves - except for the two fields that
# are genuinely new on the inserted row (Fecha / Volumen). All the other
# fields of the new row duplicate the record that is now directly below it
# (previously row 59), matching the rest of that market/category block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 59; Excel shifts rows 59-177 down to 60-178 and the
# sheet's UsedRange/dimension grows to A1:R178 automatically.
$ws.Rows.Item(59).Insert()

# Populate the newly inserted row 59 with the new weekly record.
$ws.Range('A59').Value = 7
$ws.Range('B59').Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range('C59').Value = 'Ñuble'
$ws.Range('D59').Value = 44614
$ws.Range('E59').Value = 16
$ws.Range('F59').Value = 100112017
$ws.Range('G59').Value = 'Apio'
$ws.Range('H59').Value = 'Americana (o)'
$ws.Range('I59').Value = 'Primera'
$ws.Range('J59').Value = 60
$ws.Range('K59').Value = 8000
$ws.Range('L59').Value = 8500
$ws.Range('M59').Value = 8250
$ws.Range('N59').Value = '$/docena de matas'
$ws.Range('O59').Value = 'Provincia del Elquí'
$ws.Range('P59').Value = 1375
$ws.Range('Q59').Value = 6
$ws.Range('R59').Value = 'Hortaliza'
